# Release v0.1.0-beta: Fix validation errors and update canonical URL
#
# This script updates the "Metadata" sheet (version/status/experimental/date/
# description values) and rewrites the "Include #0" sheet (renames the
# "Operation" column header to "Description", replaces the single SNOMED
# is-a/concept row with four concept rows, and drops the now-unused third
# column) of the ValueSet workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Metadata
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version: 1.0.0 -> 0.1.0
$meta.Range("B3").Value = "0.1.0"

# Status: active -> draft
$meta.Range("B6").Value = "draft"

# Experimental: previously blank -> false
$meta.Range("B7").Value = "false"

# Date: updated timestamp
$meta.Range("B8").Value = "2025-12-26T14:13:58+00:00"

# Description: previously blank -> new description text
$meta.Range("B11").Value = "Value set for nursing interventions"

# ---------------------------------------------------------------------
# Sheet 2: Include #0
# ---------------------------------------------------------------------
$inc = $wb.Worksheets.Item("Include #0")

# Clear out the previous contents (including the old column C) before
# writing the new, smaller A1:B6 layout.
$inc.Cells.Clear()

$inc.Range("A1").Value = "Concept"
$inc.Range("B1").Value = "Description"

$inc.Range("A2").Value = "71388002"
$inc.Range("B2").Value = "Procedure"

$inc.Range("A3").Value = "225358003"
$inc.Range("B3").Value = "Wound care"

$inc.Range("A4").Value = "386373004"
$inc.Range("B4").Value = "Nutrition therapy"

$inc.Range("A5").Value = ""
$inc.Range("B5").Value = ""

$inc.Range("A6").Value = "System URI"
$inc.Range("B6").Value = "http://snomed.info/sct"

# Apply the workbook's existing styles: row 1 uses the header style (same
# style as row 1 used previously), the rest use the standard body style.
$inc.Range("A1:B1").Style = $meta.Range("A1:B1").Style
$inc.Range("A2:B6").Style = $meta.Range("A2:B2").Style

$inc.Range("A1:B6").EntireColumn.AutoFit() | Out-Null
